$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Objetivos:" row - B/C long mission text replaced by docente name ---
$ws.Range("B10").Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range("C10").Value = '5840942 - Marco Aurélio Kondracki de Alcântara'

# --- Row 13: new "Programa resumido:" row (A13 was empty; B13/C13 held docente name, now "Semestral") ---
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

# --- Row 14: label becomes "Short syllabus:", B/C -> English short-syllabus text ---
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = '1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact'
$ws.Range("C14").Value = '1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact'

# --- Row 15: label becomes "Programa:", B/C -> "01/01/2020" (date-like text; must be stored as text, not a date) ---
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Formula = '="01/01/2020"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").Formula = '="01/01/2020"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# --- Row 16: label becomes "Syllabus:", B/C -> English short-syllabus text (reused) ---
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact'
$ws.Range("C16").Value = '1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact'

# --- Row 17: label becomes "Avaliação:", B/C cleared entirely (row now only has A17) ---
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17:C17").Clear()

# --- Row 18: label becomes "Método:", new B18/C18 cells with docente name ---
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range("C18").Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# --- Row 19: label becomes "Critério:" ---
$ws.Range("A19").Value = 'Critério:'

# --- Row 20: label becomes "Norma de recuperação:" ---
$ws.Range("A20").Value = 'Norma de recuperação:'

# --- Row 21: label becomes "Bibliografia:" ---
$ws.Range("A21").Value = 'Bibliografia:'

# --- Row heights that changed ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# --- Delete the old last row (former "Bibliografia" content row 22); nothing below it to shift ---
$ws.Rows.Item(22).Delete()
